# Clear the "x" values out of B9 and C9 (the shared string "x" becomes
# unreferenced and is dropped from sharedStrings.xml on save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9:C9").ClearContents()

# Move the active selection from B12 to D22.
$ws.Range("D22").Select()
